# Apply the "emission_name" column change to the Emissions sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Emissions")

# Insert a new column before column B, pushing "emission_unit" (and its
# value column) one column to the right.
$ws.Columns("B").Insert()

# Column A keeps its row-label style; just change its text.
$ws.Range("A2").Value = "CO2"

# New column B holds the emission_name header/value (default formatting).
$ws.Range("B2").Value = "CO2 emissions"
$ws.Range("B2").ClearFormats()
$ws.Range("B1").Value = "emission_name"

# Match the best-fit column widths Excel computed for the new columns.
$ws.Columns("B").ColumnWidth = 11.998697916666666
$ws.Columns("C").ColumnWidth = 11.166666666666666

$ws.Activate()
$ws.Range("F9").Select()
